# Update "想去人数" (want-to-go count) figures in column F across sheets,
# reflecting refreshed stats for the generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 13093
$ws1.Range("F6").Value  = 102
$ws1.Range("F7").Value  = 55
$ws1.Range("F9").Value  = 19
$ws1.Range("F10").Value = 13051
$ws1.Range("F11").Value = 303
$ws1.Range("F12").Value = 552
$ws1.Range("F13").Value = 8739
$ws1.Range("F14").Value = 7778
$ws1.Range("F18").Value = 134
$ws1.Range("F26").Value = 5220

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 6

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 13093
$ws4.Range("F7").Value  = 102
$ws4.Range("F8").Value  = 55
$ws4.Range("F10").Value = 19
$ws4.Range("F11").Value = 13051
$ws4.Range("F12").Value = 303
$ws4.Range("F13").Value = 552
$ws4.Range("F14").Value = 8739
$ws4.Range("F15").Value = 7778
$ws4.Range("F19").Value = 134
$ws4.Range("F23").Value = 6
$ws4.Range("F29").Value = 5220
